# Apply the Monte Carlo / NLP power-flow related input-data edits:
#  - DGs: remove the duplicate DG row (row 3)
#  - Storages: remove the second storage row (row 3)
#  - NodeData: change Smax column (E) from 500000 to 2000 for all data rows,
#              and make NodeData the active sheet/tab
#  - PVData: move the selection to F39
#  - WindGeneration: no longer the active tab (handled automatically once
#    NodeData becomes active last)

$wb = $excel.ActiveWorkbook

# --- DGs: delete row 3 (duplicate DG definition) ---
$wsDGs = $wb.Worksheets.Item("DGs")
$wsDGs.Range("A3").EntireRow.Delete() | Out-Null
$wsDGs.Activate()
$wsDGs.Range("H14").Select() | Out-Null

# --- Storages: delete row 3 (second storage definition) ---
$wsStorages = $wb.Worksheets.Item("Storages")
$wsStorages.Range("A3").EntireRow.Delete() | Out-Null
$wsStorages.Activate()
$wsStorages.Range("A3:E3").Select() | Out-Null

# --- PVData: move selection only ---
$wsPVData = $wb.Worksheets.Item("PVData")
$wsPVData.Activate()
$wsPVData.Range("F39").Select() | Out-Null

# --- NodeData: change Smax (column E) values from 500000 to 2000 ---
$wsNodeData = $wb.Worksheets.Item("NodeData")
$wsNodeData.Range("E2:E35").Value = 2000

# Make NodeData the final active sheet/tab, with the new selection
$wsNodeData.Activate()
$wsNodeData.Range("K27").Select() | Out-Null

$wb.Save()
Write-Output "edits applied"
